# Fix Critical, High, and Medium QA issues: Config errors, CPO inflation,
# Prod zones, code cleanup, imports, and logging.
#
# Rename the workbook's tabs so they use spaces instead of underscores,
# then repoint the two chart series formulas (which live on the renamed
# sheets) at the new sheet names.

$wb = $excel.ActiveWorkbook

# --- Capture the charts (and their series) before any sheets are renamed ---
$liquidityChart = $wb.Worksheets.Item("LIQUIDITY_MONITOR").ChartObjects().Item(1).Chart
$balanceChart   = $wb.Worksheets.Item("BALANCE_SHEET_HEALTH").ChartObjects().Item(1).Chart

# --- Rename the sheets: underscores -> spaces ---
$wb.Worksheets.Item("LIQUIDITY_MONITOR").Name      = "LIQUIDITY MONITOR"
$wb.Worksheets.Item("PROFIT_CONTROL").Name         = "PROFIT CONTROL"
$wb.Worksheets.Item("BALANCE_SHEET_HEALTH").Name   = "BALANCE SHEET HEALTH"
$wb.Worksheets.Item("DEBT_MANAGER").Name           = "DEBT MANAGER"
$wb.Worksheets.Item("UPLOAD_READY_FINANCE").Name   = "UPLOAD READY FINANCE"
$wb.Worksheets.Item("CROSS_REFERENCE").Name        = "CROSS REFERENCE"

# --- Update the "Liquidity Forecast" chart (lives on LIQUIDITY MONITOR) ---
$liquiditySeries = $liquidityChart.SeriesCollection().Item(1)
$liquiditySeries.Formula = "=SERIES(""Ending Cash"",'LIQUIDITY MONITOR'!`$B`$14:`$I`$14,'LIQUIDITY MONITOR'!`$B`$33:`$I`$33,1)"

# --- Update the "Solvency Gauge" chart (lives on BALANCE SHEET HEALTH) ---
$balanceSeries1 = $balanceChart.SeriesCollection().Item(1)
$balanceSeries1.Formula = "=SERIES('BALANCE SHEET HEALTH'!I12,'BALANCE SHEET HEALTH'!`$H`$13:`$H`$14,'BALANCE SHEET HEALTH'!`$I`$13:`$I`$14,1)"

$balanceSeries2 = $balanceChart.SeriesCollection().Item(2)
$balanceSeries2.Formula = "=SERIES('BALANCE SHEET HEALTH'!J12,,'BALANCE SHEET HEALTH'!`$J`$13:`$J`$14,2)"
